# DesignDocument: "use Forms" edit
# Applies the wording changes described by the commit diff.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 1) Teil 3 intro: console menu -> window with radio buttons
Replace-Text `
    "Im Konsolenprogramm kann der Benutzer in einem Menü auswählen, was er als Ergebnis haben will: " `
    "Im Fenster kann der Benutzer über Radiobuttons auswählen, was er als Ergebnis haben will: "

# 2) "Bei Menüpunkt 1" -> "Bei Auswahl 1"
Replace-Text `
    "Bei Menüpunkt 1 bekommt der Benutzer die gleiche Ausgabe wie bei Teil 2, bei der Auswahl anderer Punkte kommt die Meldung " `
    "Bei Auswahl 1 bekommt der Benutzer die gleiche Ausgabe wie bei Teil 2, bei der Auswahl anderer Punkte kommt die Meldung "

# 3) Remove the stray lastRenderedPageBreak before "Für diesen Menüpunkt..." by
#    rewriting that run's text (forces a clean run without the page-break marker).
Replace-Text `
    "Für diesen Menüpunkt ist noch keine Funktionalität vorhanden." `
    "Für diesen Menüpunkt ist noch keine Funktionalität vorhanden."

# 4) "Im Konsolenprogramm werden alle Menüpunkte..." -> "Jede Einstellung des Radiobuttons..."
Replace-Text `
    "Im Konsolenprogramm werden alle Menüpunkte mit den entsprechenden Methoden in der Lib verbunden. Jeder Menüpunkt liefert nach Eingabe entsprechende" `
    "Jede Einstellung des Radiobuttons liefert nach Eingabe entsprechende"

Replace-Text `
    " Daten zurück." `
    " Daten im Resultatfeld zurück."

# 5) Fix "Bnutzer" -> "Benutzer"
Replace-Text `
    "Der Name wird lokal gespeichert und beim nächsten Aufruf der App wird der Bnutzer mit Namen begrüßt." `
    "Der Name wird lokal gespeichert und beim nächsten Aufruf der App wird der Benutzer mit Namen begrüßt."

# 6) "Bei jeder Menüwahl und Eingabe einer PLZ" -> "Bei jedem Search-Klick"
Replace-Text `
    "Bei jeder Menüwahl und Eingabe einer PLZ werden die zuletzt eingegebenen Daten ebenfalls lokal gespeichert." `
    "Bei jedem Search-Klick werden die zuletzt eingegebenen Daten ebenfalls lokal gespeichert."

# 7) Teil 6 "Favoriten": drop the "neuer Menüpunkt" bullet entirely, and reword the PLZ list bullet.
$favRng = $d.Content.Duplicate
if ($favRng.Find.Execute("Ein neuer Menüpunkt „Favoriten“ wird hinzugefügt")) {
    $favPara = $favRng.Paragraphs(1)
    $favPara.Range.Delete() | Out-Null
}

Replace-Text `
    "Der Benutzer kann eine Liste aus Postleitzahlen anlegen, für die er bei Programmstart automatisch das Wetter angezeigt bekommt." `
    "Der Benutzer kann in einer Textbox eine Liste aus Postleitzahlen einfügen."

# 8) Insert the new bullet describing the weather-list behaviour right after the textbox bullet.
#    InsertParagraphAfter() inherits the source paragraph's formatting (style + numPr),
#    so the new bullet keeps the same NoSpacing/list numbering automatically.
$insRng = $d.Content.Duplicate
$insRng.Find.Execute("Der Benutzer kann in einer Textbox eine Liste aus Postleitzahlen einfügen.") | Out-Null
$insPara = $insRng.Paragraphs(1)
$insPara.Range.InsertParagraphAfter()
$newBulletStart = $insPara.Range.End
$newBullet = $d.Range($newBulletStart, $newBulletStart)
$newBullet.Text = "Für diese Postleitzahlen wird bei Programmstart in iner Liste automatisch das Wetter angezeigt "
